$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($ws, $addr, $val) {
    # Force numeric/percent-looking strings to remain plain text,
    # matching the inlineStr cells in the source workbook.
    $c = $ws.Range($addr)
    $c.NumberFormat = "@"
    $c.Value = $val
    $c.Style = "Normal"
}

Set-TextValue $ws "D2" "286.44"
Set-TextValue $ws "E2" "2.61%"
Set-TextValue $ws "D3" "28.76"
Set-TextValue $ws "E3" "4.80%"
Set-TextValue $ws "D4" "5.045"
Set-TextValue $ws "E4" "4.36%"
Set-TextValue $ws "D5" "0.06704"
Set-TextValue $ws "E5" "5.14%"
Set-TextValue $ws "D6" "7.332"
Set-TextValue $ws "E6" "4.40%"
$ws.Range("B7").Value = "FTXToken"
$ws.Range("C7").Value = "https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt"
Set-TextValue $ws "D7" "1.372"
Set-TextValue $ws "E7" "5.08%"
$ws.Range("B8").Value = "MXToken"
$ws.Range("C8").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
Set-TextValue $ws "D8" "0.9401"
Set-TextValue $ws "E8" "5.09%"
$ws.Range("B9").Value = "WazirX"
$ws.Range("C9").Value = "https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx"
Set-TextValue $ws "D9" "0.1560"
Set-TextValue $ws "E9" "2.66%"
$ws.Range("B10").Value = "LiechtensteinCryptoassetsExchange"
$ws.Range("C10").Value = "https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx"
Set-TextValue $ws "D10" "0.06768"
Set-TextValue $ws "E10" "14.28%"
$ws.Range("B11").Value = "MandalaExchangeToken"
$ws.Range("C11").Value = "https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx"
Set-TextValue $ws "D11" "0.07571"
Set-TextValue $ws "E11" "0.67%"
$ws.Range("B12").Value = "BitrueCoin"
$ws.Range("C12").Value = "https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr"
Set-TextValue $ws "D12" "0.02955"
Set-TextValue $ws "E12" "1.55%"
$ws.Range("B13").Value = "BitMartToken"
$ws.Range("C13").Value = "https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx"
Set-TextValue $ws "D13" "0.08996"
Set-TextValue $ws "E13" "-0.03%"
$ws.Range("B14").Value = "BitForexToken"
$ws.Range("C14").Value = "https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf"
Set-TextValue $ws "D14" "0.001590"
Set-TextValue $ws "E14" "0.92%"
$ws.Range("B15").Value = "CoinExToken"
$ws.Range("C15").Value = "https://coinranking.com/coin/APDVU0XEViZ2o+coinextoken-cet"
Set-TextValue $ws "D15" "0.04488"
Set-TextValue $ws "E15" "1.75%"
$ws.Range("B16").Value = "One"
$ws.Range("C16").Value = "https://coinranking.com/coin/6Lga5NiXX3rT+one-one"
Set-TextValue $ws "D16" "0.0006463"
Set-TextValue $ws "E16" "0.72%"
$ws.Range("B17").Value = "TigerCash"
$ws.Range("C17").Value = "https://coinranking.com/coin/6hIn06L2+tigercash-tch"
Set-TextValue $ws "D17" "0.006568"
Set-TextValue $ws "E17" "7.55%"
$ws.Range("B18").Value = "LEO"
$ws.Range("C18").Value = "https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo"
Set-TextValue $ws "D18" "3.447"
Set-TextValue $ws "E18" "-1.07%"
$ws.Range("B19").Value = "GateToken"
$ws.Range("C19").Value = "https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt"
Set-TextValue $ws "D19" "3.386"
Set-TextValue $ws "E19" "1.83%"
Set-TextValue $ws "D20" "2.246"
Set-TextValue $ws "E20" "0.94%"
Set-TextValue $ws "D21" "0.3212"
Set-TextValue $ws "E21" "2.05%"
Set-TextValue $ws "D23" "4.061"
Set-TextValue $ws "E23" "4.21%"
Set-TextValue $ws "D24" "0.1550"
Set-TextValue $ws "E24" "3.08%"
Set-TextValue $ws "D25" "0.001180"
Set-TextValue $ws "E25" "0.35%"
Set-TextValue $ws "D26" "0.004487"
Set-TextValue $ws "E26" "4.92%"
Set-TextValue $ws "D27" "0.0001244"
Set-TextValue $ws "E27" "5.53%"
Set-TextValue $ws "D28" "0.0001615"
Set-TextValue $ws "E28" "-2.32%"
Set-TextValue $ws "D40" "0.04205"
Set-TextValue $ws "E40" "3.34%"
Set-TextValue $ws "D41" "0.006709"
Set-TextValue $ws "E41" "1.09%"
Set-TextValue $ws "D42" "0.1255"
Set-TextValue $ws "E42" "-10.25%"
Set-TextValue $ws "D43" "0.002011"
Set-TextValue $ws "E43" "-4.18%"
Set-TextValue $ws "D44" "0.01230"
Set-TextValue $ws "E44" "11.84%"
Set-TextValue $ws "D45" "0.00005648"
Set-TextValue $ws "E45" "1.61%"
Set-TextValue $ws "E46" "20.74%"
Set-TextValue $ws "D47" "0.01305"
Set-TextValue $ws "E47" "-29.46%"
